$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.54411316685208
$ws.Range("C2").Value = 5.068419542088611
$ws.Range("D2").Value = 9.622403156076201
$ws.Range("F2").Value = 37.70709593666341
$ws.Range("G2").Value = 41.61080073299033
$ws.Range("H2").Value = 17.43169091620247
$ws.Range("J2").Value = 11.43319890724648
$ws.Range("K2").Value = 11.95286578398602
$ws.Range("N2").Value = 20.14647764592812

$ws.Range("B3").Value = 12.30158481872836
$ws.Range("C3").Value = 4.871746744194202
$ws.Range("D3").Value = 9.545719509806595
$ws.Range("F3").Value = 37.68850363687537
$ws.Range("G3").Value = 41.5753797004106
$ws.Range("H3").Value = 17.47102299053211
$ws.Range("J3").Value = 11.41709592714542
$ws.Range("K3").Value = 11.79085833321793
$ws.Range("N3").Value = 20.21178750672752

$ws.Range("B4").Value = 12.15340104113525
$ws.Range("C4").Value = 4.747951682780722
$ws.Range("D4").Value = 9.500381069063902
$ws.Range("F4").Value = 37.68632663181745
$ws.Range("G4").Value = 41.56551440761876
$ws.Range("H4").Value = 17.49818608876716
$ws.Range("J4").Value = 11.40953206508962
$ws.Range("K4").Value = 11.69303207586573
$ws.Range("N4").Value = 20.2537646026878

$ws.Range("B5").Value = 12.09328729776291
$ws.Range("C5").Value = 4.696823374737464
$ws.Range("D5").Value = 9.482361497423772
$ws.Range("F5").Value = 37.68776304413129
$ws.Range("G5").Value = 41.56448158848663
$ws.Range("H5").Value = 17.5100118243093
$ws.Range("J5").Value = 11.40703599513744
$ws.Range("K5").Value = 11.65363009608698
$ws.Range("N5").Value = 20.27134375043406

$ws.Range("B6").Value = 12.083324523018
$ws.Range("C6").Value = 4.688295032459756
$ws.Range("D6").Value = 9.479397403313337
$ws.Range("F6").Value = 37.68814187705778
$ws.Range("G6").Value = 41.56449044548015
$ws.Range("H6").Value = 17.51202113900326
$ws.Range("J6").Value = 11.40665698410207
$ws.Range("K6").Value = 11.64711680546661
$ws.Range("N6").Value = 20.27429136813905

$ws.Range("B7").Value = 12.15258910525166
$ws.Range("C7").Value = 4.747264789175391
$ws.Range("D7").Value = 9.500136181452442
$ws.Range("F7").Value = 37.68633659650032
$ws.Range("G7").Value = 41.56548838598685
$ws.Range("H7").Value = 17.49834251306209
$ws.Range("J7").Value = 11.40949602610718
$ws.Range("K7").Value = 11.69249874826723
$ws.Range("N7").Value = 20.25399976384602

$ws.Range("B8").Value = 12.46039091004514
$ws.Range("C8").Value = 5.001284185419268
$ws.Range("D8").Value = 9.595610930031164
$ws.Range("F8").Value = 37.698768516435
$ws.Range("G8").Value = 41.59612144667724
$ws.Range("H8").Value = 17.44462661440648
$ws.Range("J8").Value = 11.42716570856106
$ws.Range("K8").Value = 11.89669330903684
$ws.Range("N8").Value = 20.16860776410729

$ws.Range("B9").Value = 13.06594888618606
$ws.Range("C9").Value = 5.472109160266125
$ws.Range("D9").Value = 9.795892254234783
$ws.Range("F9").Value = 37.79636976874564
$ws.Range("G9").Value = 41.75040633847654
$ws.Range("H9").Value = 17.36324454640333
$ws.Range("J9").Value = 11.48014493139718
$ws.Range("K9").Value = 12.30801327228447
$ws.Range("N9").Value = 20.01598545913217

$ws.Range("B10").Value = 13.50707982201888
$ws.Range("C10").Value = 5.797658401933681
$ws.Range("D10").Value = 9.949887987711531
$ws.Range("F10").Value = 37.9125081646734
$ws.Range("G10").Value = 41.92093033637725
$ws.Range("H10").Value = 17.31811871090716
$ws.Range("J10").Value = 11.53007447440634
$ws.Range("K10").Value = 12.61393451136418
$ws.Range("N10").Value = 19.91281156750085

$ws.Range("B11").Value = 13.70595238118333
$ws.Range("C11").Value = 5.940703004367634
$ws.Range("D11").Value = 10.02118532445091
$ws.Range("F11").Value = 37.97491144667168
$ws.Range("G11").Value = 42.01080806307576
$ws.Range("H11").Value = 17.30078610065147
$ws.Range("J11").Value = 11.55513371049671
$ws.Range("K11").Value = 12.75330341157816
$ws.Range("N11").Value = 19.86780235375937

$ws.Range("B12").Value = 13.78092404067589
$ws.Range("C12").Value = 5.994100424001576
$ws.Range("D12").Value = 10.0483414599651
$ws.Range("F12").Value = 37.99990873001416
$ws.Range("G12").Value = 42.04659797900592
$ws.Range("H12").Value = 17.29468299635844
$ws.Range("J12").Value = 11.56495583180442
$ws.Range("K12").Value = 12.80605824698744
$ws.Range("N12").Value = 19.85103406873402

$ws.Range("B13").Value = 13.76479374271199
$ws.Range("C13").Value = 5.982635243657406
$ws.Range("D13").Value = 10.04248624444178
$ws.Range("F13").Value = 37.99446451622593
$ws.Range("G13").Value = 42.03881217534706
$ws.Range("H13").Value = 17.29597692041521
$ws.Range("J13").Value = 11.56282574196083
$ws.Range("K13").Value = 12.79469827844396
$ws.Range("N13").Value = 19.85463317324973

$ws.Range("B14").Value = 13.71212753262283
$ws.Range("C14").Value = 5.945111716255959
$ws.Range("D14").Value = 10.02341644552931
$ws.Range("F14").Value = 37.9769406635676
$ws.Range("G14").Value = 42.01371743642668
$ws.Range("H14").Value = 17.30027476314764
$ws.Range("J14").Value = 11.55593514185062
$ws.Range("K14").Value = 12.75764429154103
$ws.Range("N14").Value = 19.86641729712441

$ws.Range("B15").Value = 13.67982181498211
$ws.Range("C15").Value = 5.922025921104356
$ws.Range("D15").Value = 10.0117554983346
$ws.Range("F15").Value = 37.9663844376359
$ws.Range("G15").Value = 41.99857428792813
$ws.Range("H15").Value = 17.30296729565568
$ws.Range("J15").Value = 11.55175763750228
$ws.Range("K15").Value = 12.73494342683222
$ws.Range("N15").Value = 19.87367128868548

$ws.Range("B16").Value = 13.49404002687383
$ws.Range("C16").Value = 5.788204561391532
$ws.Range("D16").Value = 9.945251736464858
$ws.Range("F16").Value = 37.90862176158461
$ws.Range("G16").Value = 41.91530306883548
$ws.Range("H16").Value = 17.31931581480567
$ws.Range("J16").Value = 11.5284836169295
$ws.Range("K16").Value = 12.6048261848015
$ws.Range("N16").Value = 19.91579166911884

$ws.Range("B17").Value = 13.37955087599124
$ws.Range("C17").Value = 5.704782688157677
$ws.Range("D17").Value = 9.904757194816943
$ws.Range("F17").Value = 37.87563122570064
$ws.Range("G17").Value = 41.86736125628257
$ws.Range("H17").Value = 17.33016417326538
$ws.Range("J17").Value = 11.51480347913438
$ws.Range("K17").Value = 12.52501937250284
$ws.Range("N17").Value = 19.94212343654333

$ws.Range("B18").Value = 13.31353472566767
$ws.Range("C18").Value = 5.656327602573423
$ws.Range("D18").Value = 9.881584188687462
$ws.Range("F18").Value = 37.8575575390101
$ws.Range("G18").Value = 41.84094567669148
$ws.Range("H18").Value = 17.33670457694044
$ws.Range("J18").Value = 11.50715614970239
$ws.Range("K18").Value = 12.47913873215419
$ws.Range("N18").Value = 19.95745003767627

$ws.Range("B19").Value = 13.29115690193393
$ws.Range("C19").Value = 5.639841752625311
$ws.Range("D19").Value = 9.873759191287826
$ws.Range("F19").Value = 37.85159322122329
$ws.Range("G19").Value = 41.83220131643112
$ws.Range("H19").Value = 17.33897066466821
$ws.Range("J19").Value = 11.5046050076152
$ws.Range("K19").Value = 12.46360969418636
$ws.Range("N19").Value = 19.96267052771636

$ws.Range("B20").Value = 13.39175610045993
$ws.Range("C20").Value = 5.713712386618359
$ws.Range("D20").Value = 9.90905581083697
$ws.Range("F20").Value = 37.87904987843214
$ws.Range("G20").Value = 41.87234485794419
$ws.Range("H20").Value = 17.32897821735012
$ws.Range("J20").Value = 11.5162368998192
$ws.Range("K20").Value = 12.53351302339466
$ws.Range("N20").Value = 19.93930162269763

$ws.Range("B21").Value = 13.72760662487711
$ws.Range("C21").Value = 5.956154519827909
$ws.Range("D21").Value = 10.02901360844912
$ws.Range("F21").Value = 37.98205084381101
$ws.Range("G21").Value = 42.02104086228101
$ws.Range("H21").Value = 17.2989998813759
$ws.Range("J21").Value = 11.5579500868016
$ws.Range("K21").Value = 12.76852891486852
$ws.Range("N21").Value = 19.86294853982564

$ws.Range("B22").Value = 13.94510650568737
$ws.Range("C22").Value = 6.110097145691094
$ws.Range("D22").Value = 10.10832025429421
$ws.Range("F22").Value = 38.05732727280751
$ws.Range("G22").Value = 42.12844425057258
$ws.Range("H22").Value = 17.28209109232017
$ws.Range("J22").Value = 11.5871490477758
$ws.Range("K22").Value = 12.92197975252101
$ws.Range("N22").Value = 19.81465403638849

$ws.Range("B23").Value = 13.82923024427747
$ws.Range("C23").Value = 6.028360618033768
$ws.Range("D23").Value = 10.06591692041432
$ws.Range("F23").Value = 38.01642623001763
$ws.Range("G23").Value = 42.07019115871974
$ws.Range("H23").Value = 17.29086979529654
$ws.Range("J23").Value = 11.57138941404667
$ws.Range("K23").Value = 12.84010978519315
$ws.Range("N23").Value = 19.84028307304503

$ws.Range("B24").Value = 13.38623871925101
$ws.Range("C24").Value = 5.709676807524191
$ws.Range("D24").Value = 9.907112068980428
$ws.Range("F24").Value = 37.87750152304067
$ws.Range("G24").Value = 41.87008819823299
$ws.Range("H24").Value = 17.32951344253258
$ws.Range("J24").Value = 11.51558817229651
$ws.Range("K24").Value = 12.52967303570615
$ws.Range("N24").Value = 19.94057677863002

$ws.Range("B25").Value = 12.90244372142014
$ws.Range("C25").Value = 5.348076227992377
$ws.Range("D25").Value = 9.740425568683454
$ws.Range("F25").Value = 37.76214183451662
$ws.Range("G25").Value = 41.69859855502845
$ws.Range("H25").Value = 17.38268930164417
$ws.Range("J25").Value = 11.4638660081454
$ws.Range("K25").Value = 12.19587614869071
$ws.Range("N25").Value = 20.05569476991806
